$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: set date (A3) and time (B3) values, matching the style/number formats
# already used by row 2 (A2 = date, B2 = time).
$ws.Range("A3").NumberFormat = $ws.Range("A2").NumberFormat
$ws.Range("A3").Value = 42619

$ws.Range("B3").Value = 0.875
$ws.Range("B3").NumberFormat = $ws.Range("B2").NumberFormat

# D3: new activity text describing the kickoff of the loyalty feature
$ws.Range("D3").Value = "Ponta-pé inicial da parte de fidelidade do sistema"

# Update the active selection to D4, as in the edited workbook
$ws.Range("D4").Select()
